# Weekly refresh of Fruta / Hortaliza data:
# The rows (2-13, 15-16) get their date / volume / price / origin values
# shuffled around; row 14 is untouched. Rather than re-deriving the
# shuffle logic, apply the known target values per row directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @{ D=fecha; J=volumen; K=precio min; L=precio max; M=precio prom; O=origen; P=precio/kg }
$targets = @{
    2  = @{ D = 44497; J = 150; K = 6000;  L = 6500;  M = 6333;  O = "Región Metropolitana"; P = 253 }
    3  = @{ D = 44498; J = 220; K = 7000;  L = 7500;  M = 7273;  O = "Región Metropolitana"; P = 291 }
    4  = @{ D = 44461; J = 100; K = 13000; L = 14000; M = 13500; O = "Provincia del Elquí"; P = 540 }
    5  = @{ D = 44467; J = 100; K = 8000;  L = 9000;  M = 8500;  O = "Región Metropolitana"; P = 340 }
    6  = @{ D = 44483; J = 350; K = 5500;  L = 6000;  M = 5714;  O = "Región Metropolitana"; P = 229 }
    7  = @{ D = 44162; J = 100; K = 7500;  L = 8000;  M = 7750;  O = "Región Metropolitana"; P = 310 }
    8  = @{ D = 44335; J = 100; K = 18000; L = 20000; M = 19000; O = "Provincia de Limarí"; P = 760 }
    9  = @{ D = 44188; J = 100; K = 18000; L = 20000; M = 19000; O = "Región Metropolitana"; P = 760 }
    10 = @{ D = 44160; J = 100; K = 9000;  L = 10000; M = 9500;  O = "Región Metropolitana"; P = 380 }
    11 = @{ D = 44454; J = 100; K = 13000; L = 14000; M = 13500; O = "Provincia del Elquí"; P = 540 }
    12 = @{ D = 44316; J = 100; K = 16000; L = 18000; M = 17000; O = "Región Metropolitana"; P = 680 }
    13 = @{ D = 44351; J = 100; K = 15000; L = 16000; M = 15500; O = "Región Metropolitana"; P = 620 }
    15 = @{ D = 44482; J = 430; K = 8000;  L = 8500;  M = 8267;  O = "Región de O'Higgins"; P = 331 }
    16 = @{ D = 44384; J = 100; K = 12000; L = 13000; M = 12500; O = "Región de Coquimbo"; P = 500 }
}

foreach ($row in $targets.Keys) {
    $t = $targets[$row]
    $ws.Range("D$row").Value = $t.D
    $ws.Range("J$row").Value = $t.J
    $ws.Range("K$row").Value = $t.K
    $ws.Range("L$row").Value = $t.L
    $ws.Range("M$row").Value = $t.M
    $ws.Range("O$row").Value = $t.O
    $ws.Range("P$row").Value = $t.P
}
